$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 15782.214
$ws.Range("J43").Value = 10740
$ws.Range("L43").Value = 10740
$ws.Range("N43").Value = -10878

$ws.Range("H45").Value = 18
$ws.Range("I45").Value = 17
$ws.Range("J45").Value = 19
$ws.Range("K45").Value = 51
$ws.Range("L45").Value = 57
$ws.Range("M45").Value = 141
$ws.Range("N45").Value = -441

$ws.Range("H64").Value = 4000
$ws.Range("J64").Value = 4000
$ws.Range("L64").Value = 4000
$ws.Range("N64").Value = -4496

$ws.Range("H67").Value = 4000
$ws.Range("J67").Value = 4000
$ws.Range("L67").Value = 4000
$ws.Range("N67").Value = -5716

$ws.Range("H80").Value = 1640.5625
$ws.Range("I80").Value = 1712.375
$ws.Range("K80").Value = 5137.125
$ws.Range("M80").Value = -4139.125

$ws.Range("H83").Value = 1640.5625
$ws.Range("I83").Value = 1712.375
$ws.Range("K83").Value = 15411.375
$ws.Range("M83").Value = -10419.375

$ws.Range("H111").Value = 5351.615
$ws.Range("I111").Value = 3667.2
$ws.Range("K111").Value = 11001.6
$ws.Range("M111").Value = -7934.599999999999

$ws.Range("H113").Value = 66047.06
$ws.Range("I113").Value = 134756.25
$ws.Range("J113").Value = 4972.222
$ws.Range("K113").Value = 134756.25
$ws.Range("L113").Value = 4972.222
$ws.Range("M113").Value = -131502.25
$ws.Range("N113").Value = -11480.222

$ws.Range("H127").Value = 168052.33
$ws.Range("I127").Value = 1549.25
$ws.Range("J127").Value = 501058.5
$ws.Range("K127").Value = 4647.75
$ws.Range("L127").Value = 1503175.5
$ws.Range("M127").Value = 312.25
$ws.Range("N127").Value = -1513095.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4044.125
$ws.Range("I32").Value = 3653.7334
$ws.Range("K32").Value = 3653.7334
$ws.Range("M32").Value = -3366.7334

$ws.Range("H45").Value = 27504.4
$ws.Range("J45").Value = 23502.8
$ws.Range("L45").Value = 23502.8
$ws.Range("N45").Value = -24256.8

$ws.Range("H46").Value = 5000
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -4681
$ws.Range("N46").Value = -5638

$ws.Range("H132").Value = 3603.818
$ws.Range("I132").Value = 3488.6843
$ws.Range("K132").Value = 10466.0529
$ws.Range("M132").Value = -7936.052899999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 13884.214
$ws.Range("I86").Value = 9922.9
$ws.Range("J86").Value = 23787.5
$ws.Range("K86").Value = 9922.9
$ws.Range("L86").Value = 23787.5
$ws.Range("M86").Value = -8799.9
$ws.Range("N86").Value = -26033.5

$ws.Range("H89").Value = 13884.214
$ws.Range("I89").Value = 9922.9
$ws.Range("J89").Value = 23787.5
$ws.Range("K89").Value = 49614.5
$ws.Range("L89").Value = 118937.5
$ws.Range("M89").Value = -43998.5
$ws.Range("N89").Value = -130169.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 94355.22
$ws.Range("I31").Value = 125445.695
$ws.Range("J31").Value = 13520
$ws.Range("K31").Value = 125445.695
$ws.Range("L31").Value = 13520
$ws.Range("M31").Value = -125150.695
$ws.Range("N31").Value = -14110

$ws.Range("H34").Value = 94355.22
$ws.Range("I34").Value = 125445.695
$ws.Range("J34").Value = 13520
$ws.Range("K34").Value = 125445.695
$ws.Range("L34").Value = 13520
$ws.Range("M34").Value = -125243.695
$ws.Range("N34").Value = -13924

$ws.Range("H107").Value = 1067.4
$ws.Range("I107").Value = 1041.5555
$ws.Range("K107").Value = 1041.5555
$ws.Range("M107").Value = 878.4445000000001

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 1
$ws.Range("I20").Value = 1
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 224

$ws.Range("H93").Value = 6571.9
$ws.Range("I93").Value = 3144.8
$ws.Range("J93").Value = 9999
$ws.Range("K93").Value = 9434.400000000001
$ws.Range("L93").Value = 29997
$ws.Range("M93").Value = -7562.400000000001
$ws.Range("N93").Value = -33741

$ws.Range("H139").Value = 6142.5713
$ws.Range("I139").Value = 6499.6665
$ws.Range("J139").Value = 4000
$ws.Range("K139").Value = 19498.9995
$ws.Range("L139").Value = 12000
$ws.Range("M139").Value = -14358.9995
$ws.Range("N139").Value = -22280

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 20033
$ws.Range("J52").Value = 20033
$ws.Range("L52").Value = 20033
$ws.Range("N52").Value = -20551

$ws.Range("H80").Value = 2951
$ws.Range("J80").Value = 3001.6667
$ws.Range("L80").Value = 3001.6667
$ws.Range("N80").Value = -4997.6667

$ws.Range("H83").Value = 2951
$ws.Range("J83").Value = 3001.6667
$ws.Range("L83").Value = 15008.3335
$ws.Range("N83").Value = -24992.3335

$ws.Range("H97").Value = 1526.7142
$ws.Range("I97").Value = 1447.9166
$ws.Range("K97").Value = 1447.9166
$ws.Range("M97").Value = -951.9166

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3670.7273
$ws.Range("I40").Value = 3297.6667
$ws.Range("K40").Value = 3297.6667
$ws.Range("M40").Value = -3161.6667

$ws.Range("H46").Value = 1740
$ws.Range("I46").Value = 1740
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1740
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1552
$ws.Range("N46").ClearContents()

$ws.Range("H53").Value = 23046.8
$ws.Range("I53").Value = 10228
$ws.Range("J53").Value = 42275
$ws.Range("K53").Value = 10228
$ws.Range("L53").Value = 42275
$ws.Range("M53").Value = -9710
$ws.Range("N53").Value = -43311

$ws.Range("H122").Value = 563517.75
$ws.Range("I122").Value = 719090.8
$ws.Range("J122").Value = 19012.25
$ws.Range("K122").Value = 2157272.4
$ws.Range("L122").Value = 57036.75
$ws.Range("M122").Value = -2154822.4
$ws.Range("N122").Value = -61936.75
